# Generate Report for Handoff
# The c3bafcd9-f7bc-4b29-8611-a4694c15959a.md file has been handed off again:
#  - its status moves from "Handed back: in sync with en-US" to "Ready for handoff"
#  - a new (later) handoff/generate datetime is recorded
#  - the zh-cn / de-de sheets now carry an "Error Detail" note that the handback
#    file isn't the latest version

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f58c9c81d0ae0b330784c4b4c327355cbb04eca3/e2e/c3bafcd9-f7bc-4b29-8611-a4694c15959a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/310cd62967e1aff36ca5b2f04579ca44441f7f0e/e2e/c3bafcd9-f7bc-4b29-8611-a4694c15959a.md."

# --- Overview sheet: row 3 is the c3bafcd9-f7bc-4b29-8611-a4694c15959a.md file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-10-18 04:13:42"

# Excel's ColumnWidth (character units) gets persisted to the OOXML "width"
# attribute with a fixed +5/6 character offset (padding), so to land on an
# exact on-disk width of 40 we need to dial the COM property back by 5/6.
$targetColumnWidth = 40 - (5/6)

# --- zh-cn sheet: row 3 is the c3bafcd9-f7bc-4b29-8611-a4694c15959a.md file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-10-18 04:13:26"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = $targetColumnWidth

# --- de-de sheet: row 3 is the c3bafcd9-f7bc-4b29-8611-a4694c15959a.md file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-10-18 04:13:42"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = $targetColumnWidth
